# Update the LinkedList timing values in column C of Sheet1.
# These cells feed the embedded chart's cached series values, which
# Excel will refresh automatically to match the new cell contents.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C17").Value = 7.38
$ws.Range("C18").Value = 8.34
$ws.Range("C21").Value = 11.59

# Move the active selection, mirroring the final cursor position left
# behind after the edit (was F9, now C22).
$ws.Range("C22").Select()
